# Regenerate merged AHB files
# - rename the "_old" / "_new" header suffixes to "_FV2404" / "_FV2410"
# - turn the data range into an Excel Table ("Table1")
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels (row 1) ---------------------------------
# Columns A:J used the "_old" suffix, columns L:U used the "_new" suffix.
# Column K ("diff") stays untouched.

$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Text -replace "_old$", "_FV2404")
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Text -replace "_new$", "_FV2410")
}

# --- 2. Convert the used range into a Table ------------------------------

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U67"), [Type]::Missing, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ---------------------------------------------

[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
